$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 80

# Columns A and D hold values that look like a date ("2024-01-20") and a
# number ("02") respectively, but the source data keeps them as literal
# text (leading zero preserved, ISO date string kept verbatim). Force the
# cells to a Text number format first so Excel's autodetection doesn't
# silently convert them into a date serial / numeric value.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2024-01-20"
$ws.Cells.Item($row, 2).Value = "19:10:49"
$ws.Cells.Item($row, 3).Value = "Saturday"
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "02"
$ws.Cells.Item($row, 5).Value = 138797
$ws.Cells.Item($row, 6).Value = 140785
$ws.Cells.Item($row, 7).Value = 171734
$ws.Cells.Item($row, 8).Value = 148871
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 122941
$ws.Cells.Item($row, 11).Value = 223822
$ws.Cells.Item($row, 12).Value = 255635
$ws.Cells.Item($row, 13).Value = 185287
$ws.Cells.Item($row, 14).Value = 110373
$ws.Cells.Item($row, 15).Value = 41231
$ws.Cells.Item($row, 16).Value = 30925
$ws.Cells.Item($row, 17).Value = 73636
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42878
$ws.Cells.Item($row, 20).Value = -1
